$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.112.48"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.558.92"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'605.59"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'144.40"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").Value = "3.558.21"
$ws.Range("E7").Value = "  +3.74%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "'7.89"
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "4.166.58"
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").Value = "'29.98"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "3.558.20"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "66.222.63"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "'11.32"
$ws.Range("E19").Value = "  +9.75%  "
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("E23").Value = "  +5.60%  "
$ws.Range("D24").Value = "'79.15"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").Value = "3.704.78"
$ws.Range("E25").Value = "  +4.58%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +6.85%  "
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").Value = "'7.94"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'9.07"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'25.56"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.46"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "3.556.80"
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("D35").Value = "'0.152"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("E38").Value = "  +4.31%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "'175.68"
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("D42").Value = "'0.0848"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").Value = "'0.892"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "'46.03"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'25.74"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "'23.50"
$ws.Range("E49").Value = "  +12.28%  "
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("E51").Value = "  +0.58%  "
